# Remove the bogus "W" (NaN) book block — rows 20:25 — which also
# shrinks the used range from A1:K25 down to A1:K19 and drops the
# associated merged-cell ranges automatically.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("20:25").Delete()

# Reassign the order of the summary columns (B:E) for every surviving
# book block. Old order was Cuasidesviación, Número de Valoraciones,
# Mediana, Media; new order is Número de Valoraciones, Media,
# Cuasidesviación, Mediana.

# Header row
$ws.Range("B1").Value = "Número de Valoraciones"
$ws.Range("C1").Value = "Media"
$ws.Range("D1").Value = "Cuasidesviación"
$ws.Range("E1").Value = "Mediana"

# Libro 18 block (row 2, merged down to row 7)
$ws.Range("B2").Value = 132
$ws.Range("C2").Value = 5.727272727272728
$ws.Range("D2").Value = 0.6189243203857999
$ws.Range("E2").Value = 6

# Libro 25 block (row 8, merged down to row 13)
$ws.Range("B8").Value = 68
$ws.Range("C8").Value = 5.823529411764706
$ws.Range("D8").Value = 0.5166244188642394
$ws.Range("E8").Value = 6

# Libro 42 block (row 14, merged down to row 19)
$ws.Range("B14").Value = 156
$ws.Range("C14").Value = 5.564102564102564
$ws.Range("D14").Value = 0.7802355903888978
$ws.Range("E14").Value = 6
